$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-roster the 2v2 player names (row order unchanged, new handles)
$ws.Range("A2").Value = "Joebro"
$ws.Range("A5").Value = "MaggotEatr"
$ws.Range("A6").Value = "Sweaty"
$ws.Range("A7").Value = "ULTRA BEAST"
$ws.Range("A8").Value = "HyperDonk"
$ws.Range("A9").Value = "dverad74"
$ws.Range("A13").Value = "Pater854321"

# Match the saved selection left by the editing session
$ws.Range("B19").Select() | Out-Null
